$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 801.5
$ws.Range("I12").Value = 801.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 801.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -631.5
$ws.Range("N12").ClearContents()

$ws.Range("H92").Value = 1343.75
$ws.Range("I92").Value = 857.1667
$ws.Range("K92").Value = 857.1667
$ws.Range("M92").Value = 390.8333

$ws.Range("H123").Value = 199964.5
$ws.Range("J123").Value = 199964.5
$ws.Range("L123").Value = 199964.5
$ws.Range("N123").Value = -209764.5

$ws.Range("H132").Value = 4485.25
$ws.Range("I132").Value = 4669.6763
$ws.Range("K132").Value = 14009.0289
$ws.Range("M132").Value = -11479.0289

$ws.Range("H137").Value = 2465.8
$ws.Range("I137").Value = 2368.3845
$ws.Range("J137").Value = 3099
$ws.Range("K137").Value = 7105.1535
$ws.Range("L137").Value = 9297
$ws.Range("M137").Value = -4555.1535
$ws.Range("N137").Value = -14397

$ws.Range("H138").Value = 3599.2307
$ws.Range("J138").Value = 3262.1853
$ws.Range("L138").Value = 9786.555899999999
$ws.Range("N138").Value = -20066.5559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 285
$ws.Range("J4").Value = 270
$ws.Range("L4").Value = 270
$ws.Range("N4").Value = -502

$ws.Range("H61").Value = 3488.2917
$ws.Range("I61").Value = 2408.6428
$ws.Range("J61").Value = 4999.8
$ws.Range("K61").Value = 2408.6428
$ws.Range("L61").Value = 4999.8
$ws.Range("M61").Value = -2196.6428
$ws.Range("N61").Value = -5423.8

$ws.Range("H74").Value = 134796.22
$ws.Range("I74").Value = 200518.17
$ws.Range("J74").Value = 3352.2856
$ws.Range("K74").Value = 200518.17
$ws.Range("L74").Value = 3352.2856
$ws.Range("M74").Value = -199644.17
$ws.Range("N74").Value = -5100.2856

$ws.Range("H77").Value = 134796.22
$ws.Range("I77").Value = 200518.17
$ws.Range("J77").Value = 3352.2856
$ws.Range("K77").Value = 1002590.85
$ws.Range("L77").Value = 16761.428
$ws.Range("M77").Value = -998222.8500000001
$ws.Range("N77").Value = -25497.428

$ws.Range("H88").Value = 3571.4167
$ws.Range("I88").Value = 1747.5
$ws.Range("K88").Value = 1747.5
$ws.Range("M88").Value = -1341.5

$ws.Range("H91").Value = 3571.4167
$ws.Range("I91").Value = 1747.5
$ws.Range("K91").Value = 1747.5
$ws.Range("M91").Value = -343.5

$ws.Range("H122").Value = 2970.7058
$ws.Range("I122").Value = 2957.077
$ws.Range("J122").Value = 3015
$ws.Range("K122").Value = 8871.231
$ws.Range("L122").Value = 9045
$ws.Range("M122").Value = -6421.231
$ws.Range("N122").Value = -13945

$ws.Range("H136").Value = 3488.2917
$ws.Range("I136").Value = 2408.6428
$ws.Range("J136").Value = 4999.8
$ws.Range("K136").Value = 7225.928400000001
$ws.Range("L136").Value = 14999.4
$ws.Range("M136").Value = -4675.928400000001
$ws.Range("N136").Value = -20099.4

$ws.Range("H138").Value = 119011
$ws.Range("J138").Value = 119011
$ws.Range("L138").Value = 119011
$ws.Range("N138").Value = -129291

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 52631880
$ws.Range("J80").Value = 350.35715
$ws.Range("L80").Value = 350.35715
$ws.Range("N80").Value = -2346.35715

$ws.Range("H83").Value = 52631880
$ws.Range("J83").Value = 350.35715
$ws.Range("L83").Value = 1751.78575
$ws.Range("N83").Value = -11735.78575

$ws.Range("H86").Value = 2943.72
$ws.Range("I86").Value = 2439.2856
$ws.Range("K86").Value = 2439.2856
$ws.Range("M86").Value = -1316.2856

$ws.Range("H89").Value = 2943.72
$ws.Range("I89").Value = 2439.2856
$ws.Range("K89").Value = 12196.428
$ws.Range("M89").Value = -6580.428

$ws.Range("H99").Value = 127687.875
$ws.Range("I99").Value = 168583.83
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 168583.83
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -167085.83
$ws.Range("N99").Value = -7996

$ws.Range("H105").Value = 10834945
$ws.Range("I105").Value = 626390.2
$ws.Range("K105").Value = 626390.2
$ws.Range("M105").Value = -624643.2

$ws.Range("H107").Value = 3664315
$ws.Range("I107").Value = 6411337
$ws.Range("J107").Value = 1618.6666
$ws.Range("K107").Value = 6411337
$ws.Range("L107").Value = 1618.6666
$ws.Range("M107").Value = -6409417
$ws.Range("N107").Value = -5458.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 25001384
$ws.Range("I16").Value = 35715388
$ws.Range("K16").Value = 35715388
$ws.Range("M16").Value = -35715101

$ws.Range("H22").Value = 519.8333
$ws.Range("I22").Value = 569.75
$ws.Range("K22").Value = 569.75
$ws.Range("M22").Value = -219.75

$ws.Range("H31").Value = 6393.3076
$ws.Range("J31").Value = 9001.166999999999
$ws.Range("L31").Value = 9001.166999999999
$ws.Range("N31").Value = -9591.166999999999

$ws.Range("H34").Value = 6393.3076
$ws.Range("J34").Value = 9001.166999999999
$ws.Range("L34").Value = 9001.166999999999
$ws.Range("N34").Value = -9405.166999999999

$ws.Range("H107").Value = 1725128.5
$ws.Range("I107").Value = 2500745
$ws.Range("K107").Value = 2500745
$ws.Range("M107").Value = -2498825

$ws.Range("H113").Value = 25001384
$ws.Range("I113").Value = 35715388
$ws.Range("K113").Value = 35715388
$ws.Range("M113").Value = -35713218

$ws.Range("H132").Value = 10421444
$ws.Range("I132").Value = 2142.1304
$ws.Range("K132").Value = 6426.3912
$ws.Range("M132").Value = -3896.3912

$ws.Range("H134").Value = 2379.8386
$ws.Range("I134").Value = 2399.4285
$ws.Range("J134").Value = 2338.7
$ws.Range("K134").Value = 7198.2855
$ws.Range("L134").Value = 7016.099999999999
$ws.Range("M134").Value = -4663.2855
$ws.Range("N134").Value = -12086.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1054.8889
$ws.Range("I5").Value = 913.5714
$ws.Range("K5").Value = 2740.7142
$ws.Range("M5").Value = -2628.7142

$ws.Range("H14").Value = 471.375
$ws.Range("I14").Value = 471.375
$ws.Range("K14").Value = 1414.125
$ws.Range("M14").Value = -1241.125

$ws.Range("H98").Value = 762.25
$ws.Range("J98").Value = 899.6667
$ws.Range("L98").Value = 2699.0001
$ws.Range("N98").Value = -5695.0001

$ws.Range("H135").Value = 1054.8889
$ws.Range("I135").Value = 913.5714
$ws.Range("K135").Value = 8222.142600000001
$ws.Range("M135").Value = -5687.142600000001

$ws.Range("H140").Value = 6211.028
$ws.Range("I140").Value = 2786.5667
$ws.Range("K140").Value = 8359.7001
$ws.Range("M140").Value = -3179.7001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5509
$ws.Range("I70").Value = 4181.1816
$ws.Range("K70").Value = 4181.1816
$ws.Range("M70").Value = -3911.1816

$ws.Range("H73").Value = 5509
$ws.Range("I73").Value = 4181.1816
$ws.Range("K73").Value = 4181.1816
$ws.Range("M73").Value = -3245.1816

$ws.Range("H122").Value = 2625.6177
$ws.Range("I122").Value = 2400.818
$ws.Range("J122").Value = 3037.75
$ws.Range("K122").Value = 7202.454000000001
$ws.Range("L122").Value = 9113.25
$ws.Range("M122").Value = -4752.454000000001
$ws.Range("N122").Value = -14013.25

$ws.Range("H132").Value = 1705.2894
$ws.Range("I132").Value = 1473.3572
$ws.Range("J132").Value = 2354.7
$ws.Range("K132").Value = 4420.071599999999
$ws.Range("L132").Value = 7064.099999999999
$ws.Range("M132").Value = -1890.071599999999
$ws.Range("N132").Value = -12124.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 60201.473
$ws.Range("I40").Value = 73588.664
$ws.Range("K40").Value = 73588.664
$ws.Range("M40").Value = -73452.664

$ws.Range("H68").Value = 3267.4119
$ws.Range("I68").Value = 3229.9167
$ws.Range("K68").Value = 3229.9167
$ws.Range("M68").Value = -2480.9167

$ws.Range("H71").Value = 3267.4119
$ws.Range("I71").Value = 3229.9167
$ws.Range("K71").Value = 16149.5835
$ws.Range("M71").Value = -12405.5835

$ws.Range("H93").Value = 1162.88
$ws.Range("I93").Value = 1170.2858
$ws.Range("K93").Value = 1170.2858
$ws.Range("M93").Value = 77.71419999999989

$ws.Range("H122").Value = 10338.074
$ws.Range("I122").Value = 8034.533
$ws.Range("J122").Value = 13217.5
$ws.Range("K122").Value = 24103.599
$ws.Range("L122").Value = 39652.5
$ws.Range("M122").Value = -21653.599
$ws.Range("N122").Value = -44552.5

$ws.Range("H132").Value = 5537.643
$ws.Range("I132").Value = 3018.5
$ws.Range("J132").Value = 13598.9
$ws.Range("K132").Value = 9055.5
$ws.Range("L132").Value = 40796.7
$ws.Range("M132").Value = -6525.5
$ws.Range("N132").Value = -45856.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 166670370
$ws.Range("I100").Value = 5150.5
$ws.Range("J100").Value = 250002980
$ws.Range("K100").Value = 10301
$ws.Range("L100").Value = 500005960
$ws.Range("M100").Value = -9760
$ws.Range("N100").Value = -500007042

$ws.Range("H122").Value = 12501983
$ws.Range("I122").Value = 1939.5714
$ws.Range("K122").Value = 5818.7142
$ws.Range("M122").Value = -3368.7142
